$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Header row changes: D1 empresa -> servicio, F1 cuenta -> empresa
$ws.Range("D1").Value = "servicio"
$ws.Range("F1").Value = "empresa"

# New test case row 7 (ATC06_descargarCartolaLuz) gains a servicio value (D7) and cuenta/empresa value (F7)
$ws.Range("D7").Value = "Luz"
$ws.Range("F7").Value = "CGE"

# Update selection to D8 (matches new UI state in the diff)
$ws.Range("D8").Select()
